$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 112381898
$ws.Range("B2").Value = 90837
$ws.Range("E2").Value = 5966
$ws.Range("F2").Value = "Motaggsvamp"
$ws.Range("G2").Value = "Sarcodon squamosus"
$ws.Range("H2").Value = "(Schaeff.) Quél."
$ws.Range("Q2").Value = 608205
$ws.Range("R2").Value = 7225442

# Row 3
$ws.Range("B3").Value = 90830

# Row 4
$ws.Range("A4").Value = 112381886
$ws.Range("B4").Value = 90802
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 149
$ws.Range("F4").Value = "Tallgråticka"
$ws.Range("G4").Value = "Boletopsis grisea"
$ws.Range("H4").Value = "(Peck) Bondartsev & Singer"
$ws.Range("Q4").Value = 608475
$ws.Range("R4").Value = 7225446

# Row 5
$ws.Range("A5").Value = 112381887
$ws.Range("B5").Value = 90830
$ws.Range("E5").Value = 2059
$ws.Range("F5").Value = "Skrovlig taggsvamp"
$ws.Range("G5").Value = "Hydnellum scabrosum"
$ws.Range("H5").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q5").Value = 608478
$ws.Range("R5").Value = 7225440

# Row 6
$ws.Range("A6").Value = 112381899
$ws.Range("B6").Value = 89573
$ws.Range("E6").Value = 5442
$ws.Range("F6").Value = "Tallticka"
$ws.Range("G6").Value = "Porodaedalea pini"
$ws.Range("H6").Value = "(Brot.) Murrill"
$ws.Range("Q6").Value = 608138
$ws.Range("R6").Value = 7225544

# Row 7
$ws.Range("A7").Value = 112381884
$ws.Range("B7").Value = 77650
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("Q7").Value = 608527
$ws.Range("R7").Value = 7225507

# Row 8
$ws.Range("B8").Value = 90830
